$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header for the Percentage column
$ws.Range("F3").Value = "Percentage"

# Total column (E): first cell typed individually, then the rest filled down
# as one formula entry (mirrors how a user fills a formula down a range in Excel).
$ws.Range("E4").Formula = "=B4+C4+D4"
$ws.Range("E5:E7").Formula = "=B5+C5+D5"

# Percentage column (F): share of the grand total (E8) for each category.
$ws.Range("F4").Formula = "=E4/`$E`$8"
$ws.Range("F5:F7").Formula = "=E5/`$E`$8"

# Totals row (row 8)
$ws.Range("B8").Formula = "=B4+B5+B6+B7"
$ws.Range("C8:E8").Formula = "=C4+C5+C6+C7"
$ws.Range("E8").Formula = "=B8+C8+D8"
$ws.Range("F8").Formula = "=E8/`$E`$8"

# Matches the final selection left behind in the saved workbook
$ws.Range("H9").Select()
